# Insert a new weekly price-report row for "Poroto verde" (Macroferia
# Regional de Talca) ahead of the existing historical rows. This pushes
# the previous rows 104-106 down to 105-107 and fills the freshly
# inserted row 104 with the new week's figures (origin: Peru).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 104..106 down to 105..107, leaving a blank row 104 that
# inherits the surrounding number formatting (date style on column D).
$ws.Rows.Item(104).Insert()

$ws.Range("A104").Value = 5
$ws.Range("B104").Value = "Macroferia Regional de Talca"
$ws.Range("C104").Value = "Maule"
$ws.Range("D104").Value = 44509
$ws.Range("E104").Value = 7
$ws.Range("F104").Value = 100112031
$ws.Range("G104").Value = "Poroto verde"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 300
$ws.Range("K104").Value = 28000
$ws.Range("L104").Value = 28000
$ws.Range("M104").Value = 28000
$ws.Range("N104").Value = "$/saco 25 kilos"
$ws.Range("O104").Value = "Perú"
$ws.Range("P104").Value = 1120
$ws.Range("Q104").Value = 25
$ws.Range("R104").Value = "Hortaliza"
